$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (congruent_log/incongruent_log swapped to incongruent/congruent)
$ws.Range("D1").Value = "incongruent"
$ws.Range("E1").Value = "congruent"

# Update reaction-time log values for subsample rows 2-42

$ws.Range("D2").Value = 6.855154431973776
$ws.Range("E2").Value = 6.8632953178494
$ws.Range("D3").Value = 7.403937288194324
$ws.Range("E3").Value = 7.333570272277946
$ws.Range("D4").Value = 6.78161611301297
$ws.Range("E4").Value = 6.856756939406908
$ws.Range("D5").Value = 6.808467489811917
$ws.Range("E5").Value = 6.80157476623411
$ws.Range("D6").Value = 7.399737799218406
$ws.Range("E6").Value = 7.285273671080859
$ws.Range("D7").Value = 7.083731622688703
$ws.Range("E7").Value = 6.970550195868143
$ws.Range("D8").Value = 7.30436030485769
$ws.Range("E8").Value = 7.140148557665009
$ws.Range("D9").Value = 7.167623867124245
$ws.Range("E9").Value = 7.207124754532847
$ws.Range("D10").Value = 6.789510803817839
$ws.Range("E10").Value = 6.722531810094384
$ws.Range("D11").Value = 6.935224266524065
$ws.Range("E11").Value = 6.955017636806402
$ws.Range("D12").Value = 7.031912383117052
$ws.Range("E12").Value = 6.894791813881637
$ws.Range("D13").Value = 6.529945522631261
$ws.Range("E13").Value = 6.48114078891111
$ws.Range("D14").Value = 7.023431838841996
$ws.Range("E14").Value = 7.013676753608281
$ws.Range("D15").Value = 6.667509593992432
$ws.Range("E15").Value = 6.540272962566089
$ws.Range("D16").Value = 7.096950649446637
$ws.Range("E16").Value = 6.974750292283433
$ws.Range("D17").Value = 7.021528423761193
$ws.Range("E17").Value = 6.985212158896155
$ws.Range("D18").Value = 7.046314635041164
$ws.Range("E18").Value = 6.978568874751917
$ws.Range("D19").Value = 7.310786280138421
$ws.Range("E19").Value = 7.157297398362553
$ws.Range("D20").Value = 7.271475287597815
$ws.Range("E20").Value = 6.995286995000902
$ws.Range("D21").Value = 7.026852149233445
$ws.Range("E21").Value = 6.931727964530378
$ws.Range("D22").Value = 7.077932208860555
$ws.Range("E22").Value = 7.066079424716583
$ws.Range("D23").Value = 7.234506157161472
$ws.Range("E23").Value = 7.015097710783885
$ws.Range("D24").Value = 7.039577099527939
$ws.Range("E24").Value = 6.868211027871116
$ws.Range("D25").Value = 6.802491171612626
$ws.Range("E25").Value = 6.776345495165813
$ws.Range("D26").Value = 6.872839832394217
$ws.Range("E26").Value = 6.811074523045918
$ws.Range("D27").Value = 7.24424063913516
$ws.Range("E27").Value = 7.215446738197461
$ws.Range("D28").Value = 6.720504779346006
$ws.Range("E28").Value = 6.741366207748302
$ws.Range("D29").Value = 7.146337358395427
$ws.Range("E29").Value = 7.120270719059832
$ws.Range("D30").Value = 6.876301032019065
$ws.Range("E30").Value = 6.763521183487605
$ws.Range("D31").Value = 6.996693190904939
$ws.Range("E31").Value = 6.817423410757864
$ws.Range("D32").Value = 6.648836831487889
$ws.Range("E32").Value = 6.660101409346352
$ws.Range("D33").ClearContents()
$ws.Range("E33").Value = 6.48229724781169
$ws.Range("D34").Value = 7.026923425202499
$ws.Range("E34").Value = 6.999973651769692
$ws.Range("D35").Value = 7.082962198585097
$ws.Range("E35").Value = 6.941745336799163
$ws.Range("D36").Value = 7.048189154429294
$ws.Range("E36").Value = 6.929415257979321
$ws.Range("D37").Value = 7.293077199494363
$ws.Range("E37").Value = 7.187866913180908
$ws.Range("D38").Value = 6.96509591993144
$ws.Range("E38").Value = 6.942317143121987
$ws.Range("D39").Value = 7.429657439962671
$ws.Range("E39").Value = 7.331064601023954
$ws.Range("D40").Value = 6.977835036594626
$ws.Range("E40").Value = 6.840526099275227
$ws.Range("D41").Value = 6.902565226391424
$ws.Range("E41").Value = 6.844581679342041
$ws.Range("D42").Value = 7.078986091777253
$ws.Range("E42").Value = 7.139794371580491

# Remove the excluded subjects (rows 43-52), shrinking the sample to n = 41
$ws.Range("A43:E52").EntireRow.Delete()
